# Add a "Code" column (D) to Sheet1 that looks up the numeric school-district
# code for each riding name already present in column C, and tidy up the
# sheet view (zoom + active selection) to match the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header + data for column D ("Code")
# ---------------------------------------------------------------------------

# Numeric code for each data row (rows 2..68), aligned positionally with the
# existing column C riding values on the same row.
$codes = @(
    59015,59015,59015,59015,59005,59015,59017,59014,59017,59017,
    59006,59001,59016,59012,59033,59030,59032,59011,59025,59036,
    59035,59038,59039,59040,59034,59019,59003,59022,59008,59021,
    59042,59042,59037,59042,59037,59037,59005,59037,59029,59004,
    59024,59005,59024,59024,59041,59026,59027,59027,59029,59018,
    59018,59009,59037,59037,59013,59006,59006,59006,59010,59024,
    59028,59020,59037,59037,59028,59028,59028
)

# Copy the formatting already used by the neighbouring C / E columns so the
# new cells match the workbook's existing look (font "Lucida Grande", etc.)
# before any values are written.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

$ws.Range("C2:C68").Copy() | Out-Null
$ws.Range("D2:D68").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Header
$ws.Range("D1").Value = "Code"

# Data rows
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $codes[$i]
}

# ---------------------------------------------------------------------------
# 2. Sheet view tweaks: zoom level and active selection
# ---------------------------------------------------------------------------

$excel.ActiveWindow.Zoom = 128
$ws.Range("E6").Select()
